# Refresh cached Universalis market-data columns (H:N) across all eight
# job sheets in the Ravana_Profits workbook, per the scheduled-runner pull.
# Values are written via Range.Value; cells that no longer apply (e.g. no
# HQ listing/recipe) are cleared so the exporter omits them, matching a
# live re-pull rather than a stale zero-fill.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5
$ws.Range("H5").Value = 197.5
$ws.Range("I5").Value = 197.5
$ws.Range("K5").Value = 197.5
$ws.Range("M5").Value = -82.5

# Row 62
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""

# Row 65
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = ""

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = ""

# Row 131
$ws.Range("H131").Value = 8000
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").Value = ""

# Row 132
$ws.Range("H132").Value = 3101.2222
$ws.Range("I132").Value = 1623.5
$ws.Range("K132").Value = 4870.5
$ws.Range("M132").Value = -2340.5

# Row 135
$ws.Range("H135").Value = 3150
$ws.Range("I135").Value = 2300
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 20700
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -18165
$ws.Range("N135").Value = -41070

# Row 137
$ws.Range("H137").Value = 4498
$ws.Range("I137").Value = 749.5
$ws.Range("J137").Value = 5435.125
$ws.Range("K137").Value = 2248.5
$ws.Range("L137").Value = 16305.375
$ws.Range("M137").Value = 301.5
$ws.Range("N137").Value = -21405.375

# Row 141
$ws.Range("H141").Value = 5413
$ws.Range("I141").Value = 5413
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 16239
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -11059
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1775
$ws.Range("I45").Value = 1775
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1775
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -1398
$ws.Range("N45").Value = ""

# Row 61
$ws.Range("H61").Value = 2376.7778
$ws.Range("I61").Value = 1424.5
$ws.Range("K61").Value = 1424.5
$ws.Range("M61").Value = -1212.5

# Row 132
$ws.Range("H132").Value = 3483.6
$ws.Range("J132").Value = 3833
$ws.Range("L132").Value = 11499
$ws.Range("N132").Value = -16559

# Row 136
$ws.Range("H136").Value = 2376.7778
$ws.Range("I136").Value = 1424.5
$ws.Range("K136").Value = 4273.5
$ws.Range("M136").Value = -1723.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4349.4546
$ws.Range("I86").Value = 4330.625
$ws.Range("K86").Value = 4330.625
$ws.Range("M86").Value = -3207.625

# Row 89
$ws.Range("H89").Value = 4349.4546
$ws.Range("I89").Value = 4330.625
$ws.Range("K89").Value = 21653.125
$ws.Range("M89").Value = -16037.125

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 4055.7144
$ws.Range("I99").Value = 3878.8
$ws.Range("J99").Value = 4498
$ws.Range("K99").Value = 3878.8
$ws.Range("L99").Value = 4498
$ws.Range("M99").Value = -2380.8
$ws.Range("N99").Value = -7494

# Row 126
$ws.Range("H126").Value = 4055.7144
$ws.Range("I126").Value = 3878.8
$ws.Range("J126").Value = 4498
$ws.Range("K126").Value = 11636.4
$ws.Range("L126").Value = 13494
$ws.Range("M126").Value = -9166.400000000001
$ws.Range("N126").Value = -18434

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 90909090
$ws.Range("I55").Value = 90909090
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 272727270
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -272727093
$ws.Range("N55").Value = ""

# Row 131
$ws.Range("H131").Value = 1999
$ws.Range("J131").Value = 1999
$ws.Range("L131").Value = 5997
$ws.Range("N131").Value = -16077

# Row 136
$ws.Range("H136").Value = 3000
$ws.Range("I136").Value = 3000
$ws.Range("K136").Value = 9000
$ws.Range("M136").Value = -3900

# Row 139
$ws.Range("H139").Value = 1262.7142
$ws.Range("I139").Value = 1262.7142
$ws.Range("K139").Value = 3788.1426
$ws.Range("M139").Value = 1351.8574

# Row 140
$ws.Range("H140").Value = 2466.3333
$ws.Range("I140").Value = 2466.3333
$ws.Range("K140").Value = 7398.999899999999
$ws.Range("M140").Value = -2218.999899999999

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 5297.8335
$ws.Range("I80").Value = 5777.4
$ws.Range("J80").Value = 2900
$ws.Range("K80").Value = 5777.4
$ws.Range("L80").Value = 2900
$ws.Range("M80").Value = -4779.4
$ws.Range("N80").Value = -4896

# Row 83
$ws.Range("H83").Value = 5297.8335
$ws.Range("I83").Value = 5777.4
$ws.Range("J83").Value = 2900
$ws.Range("K83").Value = 28887
$ws.Range("L83").Value = 14500
$ws.Range("M83").Value = -23895
$ws.Range("N83").Value = -24484

# Row 126
$ws.Range("H126").Value = 1394.2
$ws.Range("I126").Value = 986.75
$ws.Range("K126").Value = 2960.25
$ws.Range("M126").Value = -490.25

# Row 132
$ws.Range("H132").Value = 2383.5293
$ws.Range("I132").Value = 1272.6666
$ws.Range("K132").Value = 3817.9998
$ws.Range("M132").Value = -1287.9998

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").Value = ""

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").Value = ""

$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 10481.833
$ws.Range("I81").Value = 7533
$ws.Range("J81").Value = 13430.667
$ws.Range("K81").Value = 15066
$ws.Range("L81").Value = 26861.334
$ws.Range("M81").Value = -14005
$ws.Range("N81").Value = -28983.334

# Row 84
$ws.Range("H84").Value = 10481.833
$ws.Range("I84").Value = 7533
$ws.Range("J84").Value = 13430.667
$ws.Range("K84").Value = 75330
$ws.Range("L84").Value = 134306.67
$ws.Range("M84").Value = -70026
$ws.Range("N84").Value = -144914.67

# Row 107
$ws.Range("H107").Value = 1034
$ws.Range("I107").Value = 999.5
$ws.Range("K107").Value = 2998.5
$ws.Range("M107").Value = -1078.5

# Row 132
$ws.Range("H132").Value = 4081.2666
$ws.Range("J132").Value = 4409.8887
$ws.Range("L132").Value = 13229.6661
$ws.Range("N132").Value = -18289.6661

# Row 136
$ws.Range("H136").Value = 3147.125
$ws.Range("I136").Value = 2472.5
$ws.Range("K136").Value = 7417.5
$ws.Range("M136").Value = -4867.5
